$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 new values (from old row 5)
$ws.Range("D2").Value = 44452
$ws.Range("J2").Value = 120

# Row 3 new values (from old row 2)
$ws.Range("D3").Value = 44453
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 2300
$ws.Range("L3").Value = 2300
$ws.Range("M3").Value = 2300
$ws.Range("P3").Value = 2300

# Row 5 new values (from old row 3)
$ws.Range("D5").Value = 44203
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 2000
$ws.Range("P5").Value = 2000
